$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.472.92"
$ws.Range("E2").Value = "  -2.39%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.991.24"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "329.12"
$ws.Range("E5").Value = "  -5.05%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.007"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4997"
$ws.Range("E7").Value = "  -4.21%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4217"
$ws.Range("E8").Value = "  -5.77%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "52.03"
$ws.Range("E9").Value = "  -3.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08911"
$ws.Range("E10").Value = "  -4.89%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.121"
$ws.Range("E11").Value = "  -5.42%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.36"
$ws.Range("E12").Value = "  -8.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.094"
$ws.Range("E13").Value = "  -6.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.974.46"
$ws.Range("E14").Value = "  -6.11%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.514"
$ws.Range("E15").Value = "  -6.62%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "96.22"
$ws.Range("E16").Value = "  -6.42%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001106"
$ws.Range("E18").Value = "  -6.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06630"
$ws.Range("E19").Value = "  -1.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.77"
$ws.Range("E20").Value = "  -8.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.006"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.963"
$ws.Range("E22").Value = "  -5.47%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.484.83"
$ws.Range("E23").Value = "  -2.37%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.88"
$ws.Range("E24").Value = "  -7.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.253"
$ws.Range("E25").Value = "  -3.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.50"
$ws.Range("E26").Value = "  -3.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.61"
$ws.Range("E27").Value = "  -7.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.527"
$ws.Range("E28").Value = "  -5.48%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.333"
$ws.Range("E29").Value = "  -8.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "128.10"
$ws.Range("E30").Value = "  -4.66%  "
$ws.Range("E31").Value = "  -9.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09947"
$ws.Range("E32").Value = "  -6.14%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.565"
$ws.Range("E33").Value = "  -12.85%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.843"
$ws.Range("E34").Value = "  -7.49%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.789"
$ws.Range("E35").Value = "  -4.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "9.581"
$ws.Range("E36").Value = "  -11.04%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02456"
$ws.Range("E37").Value = "  -7.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06343"
$ws.Range("E38").Value = "  -7.79%  "
$ws.Range("E39").Value = "  -3.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6522"
$ws.Range("E40").Value = "  -9.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.76"
$ws.Range("E41").Value = "  -8.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2068"
$ws.Range("E42").Value = "  -8.14%  "
$ws.Range("E43").Value = "  +0.01%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6352"
$ws.Range("E44").Value = "  -9.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.224"
$ws.Range("E45").Value = "  -7.44%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.36"
$ws.Range("E46").Value = "  -9.62%  "
$ws.Range("E47").Value = "  +0.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.532"
$ws.Range("E48").Value = "  -2.84%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000336"
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06988"
$ws.Range("E50").Value = "  -3.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.134"
$ws.Range("E51").Value = "  -6.99%  "
